$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shared-string fix "ИЗМ" -> "изм": the note currently lives at O30; it is
#    also being relocated to N32 in this edit (see step 3), so just retype it
#    in its new home with the corrected casing.

# 2. Clear the "1" placeholder values in row 3 (C3:J3), keep their formatting.
$ws.Range("C3:J3").ClearContents()

# 3. Row 32 scores corrected from 2 -> 5 for H32:J32; the shared SUM formula in
#    L32 and the M32 grade recompute/are updated accordingly. The "изм" marker
#    moves from O30 down to N32.
$ws.Range("O30").ClearContents()
$ws.Range("H32:J32").Value = 5
$ws.Range("M32").Value = 5
$ws.Range("N32").Value = "изм"

# 4. Clear the bottom totals row (row 33) values C33:J33 (styles stay), and
#    drop the now-unused SUM formula in L33.
$ws.Range("C33:J33").ClearContents()
$ws.Range("L33").ClearContents()

# 5. Selection moves to C33:M33 (active cell C33).
$ws.Range("C33:M33").Select() | Out-Null
